$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update PLC4 row: new IP address and port
$ws.Range("B5").Value = "192.168.96.58"
$ws.Range("C5").Value = 502

# Move the active selection to D5 as in the final saved state
$ws.Range("D5").Select()
